$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '96.664.88'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.585.13'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.50'
$ws.Range('E5').Value = '  -1.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '661.86'
$ws.Range('E6').Value = '  +1.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.57'
$ws.Range('E7').Value = '  +6.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.406'
$ws.Range('E8').Value = '  -2.92%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.06'
$ws.Range('E9').Value = '  +3.68%  '
$ws.Range('B10').Value = 'USDC'
$ws.Range('C10').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.00'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.586.08'
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.51'
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.203'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.247.68'
$ws.Range('E15').Value = '  -1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.514.78'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000257'
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.586.04'
$ws.Range('E18').Value = '  -1.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.02'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.68'
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.98'
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.499'
$ws.Range('E22').Value = '  +2.76%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '513.73'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.45'
$ws.Range('E24').Value = '  -2.55%  '
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.89'
$ws.Range('E26').Value = '  +2.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '97.26'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.84'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.777.36'
$ws.Range('E29').Value = '  -1.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.06'
$ws.Range('E30').Value = '  -2.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.149'
$ws.Range('E31').Value = '  +5.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.62'
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.182'
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +1.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.75'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.573'
$ws.Range('E37').Value = '  +1.45%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.47'
$ws.Range('E38').Value = '  +3.17%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '593.52'
$ws.Range('E39').Value = '  +5.15%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.61'
$ws.Range('E40').Value = '  +8.35%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.151'
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.87'
$ws.Range('E43').Value = '  +6.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.908'
$ws.Range('E44').Value = '  -3.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.81'
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '34.26'
$ws.Range('E47').Value = '  +4.43%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0422'
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.51'
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.61'
$ws.Range('E50').Value = '  +3.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.32'
$ws.Range('E51').Value = '  +1.33%  '
